$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 and 3 with new exam entries ---
$ws.Range("A2").Value = "CNG 242 (Quiz 1)"
$ws.Range("B2").Value = "2024-03-29 Friday 16:40"
$ws.Range("C2").Value = "I-103, I-104, I-105, I-106"

$ws.Range("A3").Value = "STAS 221 (Midterm 1)"
$ws.Range("B3").Value = "2024-03-31 Sunday 15:40"
$ws.Range("C3").Value = "S-119, S-121, S-122"

# --- Prepare rows 4-11 with the same formatting as row 3 (style index 2) ---
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new rows of data ---
$ws.Range("A4").Value = "TUR 102 (Midterm 1)"
$ws.Range("B4").Value = "2024-04-02 Tuesday 17:40"
$ws.Range("C4").Value = "SZ-22, SZ-23, SZ-24, SZ-25"

$ws.Range("A5").Value = "CNG 242 (Midterm 1)"
$ws.Range("B5").Value = "2024-04-16 Tuesday 17:40"
$ws.Range("C5").Value = "SZ-22, SZ-23, SZ-24, SZ-25"

$ws.Range("A6").Value = "CNG 232 (Midterm 1)"
$ws.Range("B6").Value = "2024-04-21 Sunday 14:40"
$ws.Range("C6").Value = "S-121, S-122, S-123"

$ws.Range("A7").Value = "CNG 280 (Midterm 1)"
$ws.Range("B7").Value = "2024-04-27 Saturday 13:40"
$ws.Range("C7").Value = "I-103, I-104, I-105, I-106, IZ-04"

$ws.Range("A8").Value = "CNG 242 (Quiz 2)"
$ws.Range("B8").Value = "2024-05-03 Friday 16:40"
$ws.Range("C8").Value = "I-103, I-104, I-105, I-106"

$ws.Range("A9").Value = "STAS 221 (Midterm 2)"
$ws.Range("B9").Value = "2024-05-12 Sunday 15:40"
$ws.Range("C9").Value = "S-119, S-121, S-122"

$ws.Range("A10").Value = "CNG 242 (Quiz 3)"
$ws.Range("B10").Value = "2024-05-17 Friday 16:40"
$ws.Range("C10").Value = "I-103, I-104, I-105, I-106"

$ws.Range("A11").Value = "CNG 242 (Quiz 4)"
$ws.Range("B11").Value = "2024-05-31 Friday 16:40"
$ws.Range("C11").Value = "I-103, I-104, I-105, I-106"

$wb.Save()
